# Insert a new weekly price-report row at row 44 (pushing the existing
# rows 44-101 down to 45-102) and populate it with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(44).Insert()

$ws.Range("A44").Value = 7
$ws.Range("B44").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C44").Value = "Ñuble"
$ws.Range("D44").Value = 44664
$ws.Range("E44").Value = 16
$ws.Range("F44").Value = 100112030
$ws.Range("G44").Value = "Poroto granado"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 120
$ws.Range("K44").Value = 19000
$ws.Range("L44").Value = 20000
$ws.Range("M44").Value = 19500
$ws.Range("N44").Value = "$/saco 25 kilos"
$ws.Range("O44").Value = "Provincia de Diguillín"
$ws.Range("P44").Value = 780
$ws.Range("Q44").Value = 25
$ws.Range("R44").Value = "Hortaliza"
